# 25. Clustering: Perform K-means clustering on the dataset.
# Add the K-Nearest Neighbors and Support Vector Regression result rows
# to the metrics table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    # Force the cell to store the value as text (matches the shared-string
    # cells used for the other numeric-looking metric values already in
    # the sheet), then reset the cell style back to the default "Normal"
    # so no extra formatting is left behind on the cell.
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 10: K-Nearest Neighbors
Set-TextValue "A10" "K-Nearest Neighbors"
Set-TextValue "B10" "251.64652259473945"
Set-TextValue "C10" "368727.3813827998"
Set-TextValue "D10" "607.2292659142837"

# Row 11: Support Vector Regression
Set-TextValue "A11" "Support Vector Regression"
Set-TextValue "B11" "11419.047297855577"
Set-TextValue "C11" "384427929.8857064"
Set-TextValue "D11" "19606.83375473221"

# Match the author's final selection state.
$ws.Range("B11").Select()
